$wb = $excel.ActiveWorkbook

$wsMembers3 = $wb.Worksheets.Item("Members_O12345678901234_3")
$wsMembers4 = $wb.Worksheets.Item("Members_O12345678901234_4")

# Fix invalid RELATION values on the "Members_..._3" sheet: the relation
# column (D) referenced bare code values; they must be prefixed with
# "code:" to be valid references.
$wsMembers3.Range("D3").Value = "code:testcode01"
$wsMembers3.Range("D4").Value = "code:testcode01"
$wsMembers3.Range("D6").Value = "code:testcode04"
$wsMembers3.Range("D7").Value = "code:testcode05"

# Same fix on the "Members_..._4" sheet, relation column is F here.
$wsMembers4.Range("F3").Value = "code:testcode08"
$wsMembers4.Range("F4").Value = "code:testcode09"
$wsMembers4.Range("F7").Value = "code:testcode12"
$wsMembers4.Range("F8").Value = "code:testcode13"
$wsMembers4.Range("F9").Value = "code:testcode13"

# Re-create the navigation / selection state recorded in the workbook:
# the user was on "Members_..._4" (selection B11), moved the selection to
# I22, then switched to "Members_..._3" and selected C18, which becomes
# the final active sheet/selection.
[void]$wsMembers4.Activate()
$wsMembers4.Range("I22").Select() | Out-Null

[void]$wsMembers3.Activate()
$wsMembers3.Range("C18").Select() | Out-Null
